# The document contains three occurrences of an "<id>...</id>" tag that
# was previously split across multiple runs (e.g. "<id>", "p012v", "_1",
# "</id>"). This collapses each occurrence back into a single run/word
# containing the full text, which is what Word naturally does when you
# find the full visible string and replace it in place (the matched
# range's runs get merged into one run carrying the formatting of the
# first run in the match).

$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p012v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p012v_1</id>", 2)
$d.Content.Find.Execute("<id>p012v_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p012v_2</id>", 2)
$d.Content.Find.Execute("<id>p012v_3</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p012v_3</id>", 2)
